$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 39: RenderToken -> TheSandbox
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.4963"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.39%  "

# Row 40: TheSandbox -> RenderToken
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.211"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.72%  "

# Row 2
$ws.Range("D2").Value = "26.501.08"
$ws.Range("E2").Value = "  +2.03%  "

# Row 3
$ws.Range("D3").Value = "1.849.26"
$ws.Range("E3").Value = "  +0.84%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.56%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "256.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -7.87%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.004"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.52%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5254"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.86%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3288"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.82%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06731"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.29%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.44"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.58%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7746"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.88%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07704"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.28%  "

# Row 13
$ws.Range("D13").Value = "1.828.29"
$ws.Range("E13").Value = "  -0.09%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.68"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.50%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.050"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.45%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.004"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.57%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.18"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.16%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.003"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.51%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007908"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.95%  "

# Row 20
$ws.Range("D20").Value = "26.561.99"
$ws.Range("E20").Value = "  +2.13%  "

# Row 21
$ws.Range("D21").Value = "2.082.14"
$ws.Range("E21").Value = "  +0.91%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.602"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.19%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.711"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.23%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.990"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.49%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.359"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.66%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "144.73"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.46%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.654"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.79%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.73%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "111.34"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.65%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.210"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.40%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.205"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.88%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08791"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.09%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04885"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.51%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.141"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.77%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.867"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.02%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7073"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.95%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.115"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.58%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01814"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.93%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "114.78"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.16%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8987"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.12%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.089"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.69%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.004"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.59%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.799"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.53%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4295"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.74%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1293"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.16%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.148"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.69%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05930"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.21%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "35.43"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.95%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "60.03"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.78%  "
